# Applies the attendance-sheet edits for the "August" sheet (DPG Operation
# Research attendance workbook): updates the Total-Classes "G" column entry
# for the summary row, fills in the previously-blank daily "G" (Saturday?)
# attendance-count cells for the eight student rows, and moves the active
# cell/selection to reflect where the editor finished working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("August")

# Total classes conducted so far (row 10, column G) bumped from 28 to 29.
$ws.Range("G10").Value = 29

# Per-student class counts in column G (previously blank cells that now
# carry an explicit 0/3 count, matching the other day columns D:F).
$ws.Range("G14").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("G16").Value = 3
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("G21").Value = 0

# Leave the selection where the author last clicked before saving.
$ws.Range("I21").Select()
